# Update New Orleans shard workbook:
#  - "hotel_info" (currently sheet index 1) becomes "review_info": all prior
#    data is cleared and replaced with just the review_info header row.
#  - "review_info" (currently sheet index 2) becomes "hotel_info": the header
#    row gains a new "State" column (after Hotel_Name) and the one data row
#    is rewritten to match, including the new "Louisiana" value.

$wb = $excel.ActiveWorkbook

# Grab the two worksheets by their *current* names before anything is renamed.
$hotelSheet = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# --- Turn the old "hotel_info" sheet into the new "review_info" sheet ---
# (renamed via a temporary name first so the two renames below never collide
# with each other, regardless of statement ordering/batching)
$hotelSheet.Cells.Clear()
$hotelSheet.Name = "__tmp_swap_name__"

$reviewHeaders = @(
    "STR",
    "reviewer_ID",
    "reviewer_name",
    "Review_ID",
    "Date_of_scraping",
    "ReviewURL",
    "Tripadvisor_gcode",
    "Tripadvisor_dcode",
    "Tripadvisor_rcode",
    "review_date",
    "review_title",
    "review_content",
    "review_rating",
    "trip_month",
    "trip_purpose",
    "value",
    "rooms",
    "Location",
    "Cleanliness",
    "Sleep Quality",
    "Service",
    "Picture(yes=1)",
    "respondent",
    "response_date",
    "response_text"
)

for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $hotelSheet.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# --- Turn the old "review_info" sheet into the new "hotel_info" sheet ---
$reviewSheet.Cells.Clear()
$reviewSheet.Name = "hotel_info"
$hotelSheet.Name = "review_info"

$hotelHeaders = @(
    "STR",
    "Hotel_Name",
    "State",
    "City",
    "Zip",
    "TA_ReviewURL",
    "Tripadvisor_Hotel_Name",
    "English_Reviews_num",
    "Local_Rank",
    "Total_Reviews_num"
)

for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $reviewSheet.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

$hotelRow = @(
    29674,
    "Fairfield Inn & Suites New Orleans Downtown French Quarter Area",
    "Louisiana",
    "New Orleans",
    70112,
    "https://www.tripadvisor.com/Hotel_Review-g60864-d93237-Reviews-Fairfield_Inn_Suites_New_Orleans_Downtown_French_Quarter_Area-New_Orleans_Louisiana.html",
    "Fairfield Inn & Suites New Orleans Downtown/French Quarter Area",
    "131",
    "77",
    "134"
)

# English_Reviews_num / Local_Rank / Total_Reviews_num (H2:J2) are stored as
# text in the source data (not numbers), so force those three cells to Text
# format before assigning, otherwise Excel auto-coerces numeric-looking
# strings back into numbers.
$reviewSheet.Range("H2:J2").NumberFormat = "@"

for ($i = 0; $i -lt $hotelRow.Length; $i++) {
    $reviewSheet.Cells.Item(2, $i + 1).Value = $hotelRow[$i]
}

# Both sheets keep their original tab positions (1st / 2nd) - only their
# names and contents were swapped above - so sheet order already matches
# the source (review_info first, hotel_info second).
